# Updates the "cryptos" price/volume table with refreshed figures, and swaps
# the NEARProtocol / ImmutableX rows (31/32) to reflect their new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to remain plain text even if the value looks numeric,
    # preserving exact formatting (e.g. trailing zeros, double-dot prices).
    $origStyle = $ws.Range($cellRef).Style
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = $origStyle
}

$ws.Range("D2").Value = '69.001.63'
$ws.Range("E2").Value = '  +2.01%  '
$ws.Range("D3").Value = '3.736.47'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  +0.14%  '
Set-TextValue "D5" '601.85'
$ws.Range("E5").Value = '  +1.39%  '
Set-TextValue "D6" '167.91'
$ws.Range("E6").Value = '  -2.17%  '
$ws.Range("D7").Value = '3.735.62'
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +3.07%  '
$ws.Range("E10").Value = '  +4.87%  '
Set-TextValue "D11" '6.32'
$ws.Range("E11").Value = '  +2.70%  '
$ws.Range("E12").Value = '  +0.35%  '
Set-TextValue "D13" '38.17'
$ws.Range("E13").Value = '  +1.92%  '
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").Value = '4.363.03'
$ws.Range("D16").Value = '3.736.15'
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").Value = '68.940.32'
$ws.Range("E17").Value = '  +2.00%  '
Set-TextValue "D18" '7.23'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("E19").Value = '  +0.88%  '
Set-TextValue "D20" '17.20'
$ws.Range("E20").Value = '  +6.65%  '
Set-TextValue "D21" '497.74'
$ws.Range("E21").Value = '  +1.98%  '
Set-TextValue "D22" '10.17'
$ws.Range("E22").Value = '  +13.25%  '
$ws.Range("E23").Value = '  +1.29%  '
Set-TextValue "D24" '85.35'
$ws.Range("E25").Value = '  -1.34%  '
Set-TextValue "D26" '0.0000140'
$ws.Range("E26").Value = '  -0.78%  '
Set-TextValue "D27" '12.30'
$ws.Range("E27").Value = '  +1.42%  '
Set-TextValue "D28" '10.15'
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  +0.57%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D31" '7.98'
$ws.Range("E31").Value = '  +4.32%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D32" '2.42'
$ws.Range("E32").Value = '  +2.02%  '
Set-TextValue "D33" '31.72'
$ws.Range("E33").Value = '  -1.77%  '
$ws.Range("D34").Value = '3.890.23'
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("E35").Value = '  +0.61%  '
$ws.Range("D36").Value = '3.669.82'
$ws.Range("E36").Value = '  -0.03%  '
Set-TextValue "D37" '0.999'
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("E39").Value = '  +2.07%  '
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("E41").Value = '  +0.51%  '
Set-TextValue "D42" '434.67'
$ws.Range("E42").Value = '  -3.32%  '
Set-TextValue "D43" '48.97'
$ws.Range("E43").Value = '  +0.31%  '
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("E45").Value = '  +1.06%  '
$ws.Range("E46").Value = '  +1.96%  '
Set-TextValue "D48" '40.51'
$ws.Range("E48").Value = '  -1.67%  '
Set-TextValue "D49" '141.99'
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("E50").Value = '  +1.82%  '
$ws.Range("D51").Value = '2.745.21'
$ws.Range("E51").Value = '  -1.44%  '
